# "updates for wed lecture"
#
# Rehearsed/recorded slide timings were applied to the deck: every slide's
# transition now auto-advances after a fixed number of milliseconds
# (p:transition/@advTm, surfaced on the object model as
# SlideShowTransition.AdvanceOnTime / .AdvanceTime), and the two slides that
# carry handwritten ink annotations (slides 3 and 4) each picked up a
# "TIMING" custom tag recording the recorded lecture position.

$p = $ppt.ActivePresentation

# Slide 1 -> advance automatically after 23.438s
$s1 = $p.Slides.Item(1)
$s1.SlideShowTransition.AdvanceOnTime = $true
$s1.SlideShowTransition.AdvanceTime = 23.438

# Slide 2 -> advance automatically after 77.959s
$s2 = $p.Slides.Item(2)
$s2.SlideShowTransition.AdvanceOnTime = $true
$s2.SlideShowTransition.AdvanceTime = 77.959

# Slide 3 -> advance automatically after 73.129s, plus a recorded TIMING tag
$s3 = $p.Slides.Item(3)
$s3.SlideShowTransition.AdvanceOnTime = $true
$s3.SlideShowTransition.AdvanceTime = 73.129
$s3.Tags.Add("TIMING", "|26")

# Slide 4 -> advance automatically after 225.110s, plus a recorded TIMING tag
$s4 = $p.Slides.Item(4)
$s4.SlideShowTransition.AdvanceOnTime = $true
$s4.SlideShowTransition.AdvanceTime = 225.110
$s4.Tags.Add("TIMING", "|5.1")
